# Update Sheets via scheduled runner: refresh market-board price snapshots
# for the Ixion_Profits leve-crafting profitability sheets (H:N columns).
$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 43
$ws.Range("H43").Value = 2250
$ws.Range("J43").Value = 2250
$ws.Range("L43").Value = 2250
$ws.Range("N43").Value = -2388

# Row 53
$ws.Range("H53").Value = 680.1429000000001
$ws.Range("I53").Value = 1574.5714
$ws.Range("J53").Value = 232.92857
$ws.Range("K53").Value = 1574.5714
$ws.Range("L53").Value = 232.92857
$ws.Range("M53").Value = -937.5714
$ws.Range("N53").Value = -1506.92857

# Row 55
$ws.Range("H55").Value = 503.86667
$ws.Range("I55").Value = 626.1818
$ws.Range("J55").Value = 167.5
$ws.Range("K55").Value = 626.1818
$ws.Range("L55").Value = 167.5
$ws.Range("M55").Value = -412.1818
$ws.Range("N55").Value = -595.5

# Row 132
$ws.Range("H132").Value = 1164.579
$ws.Range("I132").Value = 1187.6111
$ws.Range("K132").Value = 3562.8333
$ws.Range("M132").Value = -1032.8333

# Row 138
$ws.Range("H138").Value = 2569.9758
$ws.Range("I138").Value = 1033.5264
$ws.Range("J138").Value = 3867.422
$ws.Range("K138").Value = 3100.5792
$ws.Range("L138").Value = 11602.266
$ws.Range("M138").Value = 2039.4208
$ws.Range("N138").Value = -21882.266


# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4209.485
$ws.Range("I32").Value = 3486.64
$ws.Range("J32").Value = 6468.375
$ws.Range("K32").Value = 3486.64
$ws.Range("L32").Value = 6468.375
$ws.Range("M32").Value = -3199.64
$ws.Range("N32").Value = -7042.375

# Row 61
$ws.Range("H61").Value = 5915.6924
$ws.Range("I61").Value = 2375.9048
$ws.Range("J61").Value = 20782.8
$ws.Range("K61").Value = 2375.9048
$ws.Range("L61").Value = 20782.8
$ws.Range("M61").Value = -2163.9048
$ws.Range("N61").Value = -21206.8

# Row 74
$ws.Range("H74").Value = 2053.125
$ws.Range("I74").Value = 1869.4166
$ws.Range("K74").Value = 1869.4166
$ws.Range("M74").Value = -995.4166

# Row 77
$ws.Range("H77").Value = 2053.125
$ws.Range("I77").Value = 1869.4166
$ws.Range("K77").Value = 9347.083000000001
$ws.Range("M77").Value = -4979.083000000001

# Row 88
$ws.Range("H88").Value = 142859700
$ws.Range("I88").Value = 2993
$ws.Range("J88").Value = 200002380
$ws.Range("K88").Value = 2993
$ws.Range("L88").Value = 200002380
$ws.Range("M88").Value = -2587
$ws.Range("N88").Value = -200003192

# Row 91
$ws.Range("H91").Value = 142859700
$ws.Range("I91").Value = 2993
$ws.Range("J91").Value = 200002380
$ws.Range("K91").Value = 2993
$ws.Range("L91").Value = 200002380
$ws.Range("M91").Value = -1589
$ws.Range("N91").Value = -200005188

# Row 126
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

# Row 128
$ws.Range("H128").Value = 34600
$ws.Range("J128").Value = 34600
$ws.Range("L128").Value = 34600
$ws.Range("N128").Value = -44560

# Row 132
$ws.Range("H132").Value = 2889.02
$ws.Range("I132").Value = 1905.4073
$ws.Range("J132").Value = 4043.6956
$ws.Range("K132").Value = 5716.2219
$ws.Range("L132").Value = 12131.0868
$ws.Range("M132").Value = -3186.2219
$ws.Range("N132").Value = -17191.0868

# Row 136
$ws.Range("H136").Value = 5915.6924
$ws.Range("I136").Value = 2375.9048
$ws.Range("J136").Value = 20782.8
$ws.Range("K136").Value = 7127.714399999999
$ws.Range("L136").Value = 62348.39999999999
$ws.Range("M136").Value = -4577.714399999999
$ws.Range("N136").Value = -67448.39999999999


# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 2680.4
$ws.Range("I86").Value = 2334
$ws.Range("J86").Value = 3200
$ws.Range("K86").Value = 2334
$ws.Range("L86").Value = 3200
$ws.Range("M86").Value = -1211
$ws.Range("N86").Value = -5446

# Row 89
$ws.Range("H89").Value = 2680.4
$ws.Range("I89").Value = 2334
$ws.Range("J89").Value = 3200
$ws.Range("K89").Value = 11670
$ws.Range("L89").Value = 16000
$ws.Range("M89").Value = -6054
$ws.Range("N89").Value = -27232

# Row 134
$ws.Range("H134").Value = 2225.9062
$ws.Range("I134").Value = 2095.1875
$ws.Range("J134").Value = 2356.625
$ws.Range("K134").Value = 6285.5625
$ws.Range("L134").Value = 7069.875
$ws.Range("M134").Value = -3750.5625
$ws.Range("N134").Value = -12139.875


# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 132
$ws.Range("H132").Value = 2182.4866
$ws.Range("I132").Value = 1291.5333
$ws.Range("J132").Value = 6000.857
$ws.Range("K132").Value = 3874.5999
$ws.Range("L132").Value = 18002.571
$ws.Range("M132").Value = -1344.5999
$ws.Range("N132").Value = -23062.571


# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 17
$ws.Range("H17").Value = 1000000
$ws.Range("I17").Value = 1000000
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 3000000
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()
$ws.Range("M17").Value = -2999831

# Row 23
$ws.Range("H23").Value = 12500124
$ws.Range("I23").Value = 33333376
$ws.Range("J23").Value = 172.2
$ws.Range("K23").Value = 100000128
$ws.Range("L23").Value = 516.5999999999999
$ws.Range("M23").Value = -99999893
$ws.Range("N23").Value = -986.5999999999999

# Row 34
$ws.Range("H34").Value = 2178.6191
$ws.Range("I34").Value = 760
$ws.Range("K34").Value = 2280
$ws.Range("M34").Value = -2196

# Row 55
$ws.Range("H55").Value = 4575.3125
$ws.Range("J55").Value = 5157.5
$ws.Range("L55").Value = 15472.5
$ws.Range("N55").Value = -15826.5

# Row 68
$ws.Range("H68").Value = 3529.0942
$ws.Range("I68").Value = 4655.2593
$ws.Range("J68").Value = 2359.6155
$ws.Range("K68").Value = 13965.7779
$ws.Range("L68").Value = 7078.8465
$ws.Range("M68").Value = -13154.7779
$ws.Range("N68").Value = -8700.8465

# Row 71
$ws.Range("H71").Value = 3529.0942
$ws.Range("I71").Value = 4655.2593
$ws.Range("J71").Value = 2359.6155
$ws.Range("K71").Value = 41897.3337
$ws.Range("L71").Value = 21236.5395
$ws.Range("M71").Value = -37841.3337
$ws.Range("N71").Value = -29348.5395

# Row 107
$ws.Range("H107").Value = 750.16437
$ws.Range("I107").Value = 265.36
$ws.Range("J107").Value = 1002.6667
$ws.Range("K107").Value = 796.08
$ws.Range("L107").Value = 3008.0001
$ws.Range("M107").Value = 1123.92
$ws.Range("N107").Value = -6848.0001

# Row 113
$ws.Range("H113").Value = 270771.78
$ws.Range("I113").Value = 482.92
$ws.Range("K113").Value = 1448.76
$ws.Range("M113").Value = 721.24

# Row 132
$ws.Range("H132").Value = 2023018.1
$ws.Range("I132").Value = 1817.3334
$ws.Range("J132").Value = 2270512
$ws.Range("K132").Value = 16356.0006
$ws.Range("L132").Value = 20434608
$ws.Range("M132").Value = -13826.0006
$ws.Range("N132").Value = -20439668


# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 6050.6924
$ws.Range("I70").Value = 6062.4062
$ws.Range("J70").Value = 5997.143
$ws.Range("K70").Value = 6062.4062
$ws.Range("L70").Value = 5997.143
$ws.Range("M70").Value = -5792.4062
$ws.Range("N70").Value = -6537.143

# Row 73
$ws.Range("H73").Value = 6050.6924
$ws.Range("I73").Value = 6062.4062
$ws.Range("J73").Value = 5997.143
$ws.Range("K73").Value = 6062.4062
$ws.Range("L73").Value = 5997.143
$ws.Range("M73").Value = -5126.4062
$ws.Range("N73").Value = -7869.143


# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 20834836
$ws.Range("I46").Value = 37037940
$ws.Range("J46").Value = 2271.4285
$ws.Range("K46").Value = 37037940
$ws.Range("L46").Value = 2271.4285
$ws.Range("M46").Value = -37037752
$ws.Range("N46").Value = -2647.4285

# Row 132
$ws.Range("H132").Value = 10836138
$ws.Range("I132").Value = 13980969
$ws.Range("J132").Value = 3944.2222
$ws.Range("K132").Value = 41942907
$ws.Range("L132").Value = 11832.6666
$ws.Range("M132").Value = -41940377
$ws.Range("N132").Value = -16892.6666

# Row 136
$ws.Range("H136").Value = 4604.8423
$ws.Range("I136").Value = 2441.75
$ws.Range("K136").Value = 7325.25
$ws.Range("M136").Value = -4775.25


# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 69
$ws.Range("H69").Value = 21000
$ws.Range("I69").Value = 3000
$ws.Range("J69").Value = 30000
$ws.Range("K69").Value = 3000
$ws.Range("L69").Value = 30000
$ws.Range("M69").Value = -2251
$ws.Range("N69").Value = -31498

# Row 72
$ws.Range("H72").Value = 21000
$ws.Range("I72").Value = 3000
$ws.Range("J72").Value = 30000
$ws.Range("K72").Value = 9000
$ws.Range("L72").Value = 90000
$ws.Range("M72").Value = -5256
$ws.Range("N72").Value = -97488

# Row 136
$ws.Range("H136").Value = 2987.5312
$ws.Range("I136").Value = 3607.3333
$ws.Range("J136").Value = 2190.6428
$ws.Range("K136").Value = 10821.9999
$ws.Range("L136").Value = 6571.928400000001
$ws.Range("M136").Value = -8271.999899999999
$ws.Range("N136").Value = -11671.9284

